# Scheduled data refresh: update Leve profitability figures (columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with the latest
# market-board snapshot. Some rows gain/lose their LeveProfitNQ (M) or
# LeveProfitHQ (N) cell depending on whether that figure is computable
# with the refreshed prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11187.333
$ws.Range("I51").Value = 5685.3335
$ws.Range("J51").Value = 12287.733
$ws.Range("K51").Value = 5685.3335
$ws.Range("L51").Value = 12287.733
$ws.Range("M51").Value = -5201.3335
$ws.Range("N51").Value = -13255.733

$ws.Range("H53").Value = 241.95238
$ws.Range("I53").Value = 201.61539
$ws.Range("J53").Value = 307.5
$ws.Range("K53").Value = 201.61539
$ws.Range("L53").Value = 307.5
$ws.Range("M53").Value = 435.38461
$ws.Range("N53").Value = -1581.5

$ws.Range("H64").Value = 5500
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996

$ws.Range("H67").Value = 5500
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216

$ws.Range("H138").Value = 4177.2607
$ws.Range("J138").Value = 4460.4414
$ws.Range("L138").Value = 13381.3242
$ws.Range("N138").Value = -23661.3242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 271.08334
$ws.Range("I2").Value = 271.08334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 271.08334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -158.08334
$ws.Range("N2").ClearContents()

$ws.Range("H25").Value = 300
$ws.Range("I25").Value = 300
$ws.Range("K25").Value = 300
$ws.Range("M25").Value = 102

$ws.Range("H61").Value = 4483.6
$ws.Range("I61").Value = 1962
$ws.Range("J61").Value = 6164.6665
$ws.Range("K61").Value = 1962
$ws.Range("L61").Value = 6164.6665
$ws.Range("M61").Value = -1750
$ws.Range("N61").Value = -6588.6665

$ws.Range("H74").Value = 2674.9412
$ws.Range("I74").Value = 2030.75
$ws.Range("J74").Value = 4221
$ws.Range("K74").Value = 2030.75
$ws.Range("L74").Value = 4221
$ws.Range("M74").Value = -1156.75
$ws.Range("N74").Value = -5969

$ws.Range("H77").Value = 2674.9412
$ws.Range("I77").Value = 2030.75
$ws.Range("J77").Value = 4221
$ws.Range("K77").Value = 10153.75
$ws.Range("L77").Value = 21105
$ws.Range("M77").Value = -5785.75
$ws.Range("N77").Value = -29841

$ws.Range("H116").Value = 271.08334
$ws.Range("I116").Value = 271.08334
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 271.08334
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 2022.91666
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 1406.4445
$ws.Range("I132").Value = 1471.7646
$ws.Range("K132").Value = 4415.293799999999
$ws.Range("M132").Value = -1885.293799999999

$ws.Range("H136").Value = 4483.6
$ws.Range("I136").Value = 1962
$ws.Range("J136").Value = 6164.6665
$ws.Range("K136").Value = 5886
$ws.Range("L136").Value = 18493.9995
$ws.Range("M136").Value = -3336
$ws.Range("N136").Value = -23593.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 271.08334
$ws.Range("I3").Value = 271.08334
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 271.08334
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -157.08334
$ws.Range("N3").ClearContents()

$ws.Range("H20").Value = 600
$ws.Range("I20").Value = 600
$ws.Range("K20").Value = 600
$ws.Range("M20").Value = -353

$ws.Range("H22").Value = 613.7778
$ws.Range("I22").Value = 512.1429000000001
$ws.Range("K22").Value = 512.1429000000001
$ws.Range("M22").Value = -339.1429000000001

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18695.666
$ws.Range("I60").Value = 9111
$ws.Range("J60").Value = 23488
$ws.Range("K60").Value = 9111
$ws.Range("L60").Value = 23488
$ws.Range("M60").Value = -8600
$ws.Range("N60").Value = -24510

$ws.Range("H134").Value = 1680.8
$ws.Range("I134").Value = 1610.2222
$ws.Range("K134").Value = 4830.6666
$ws.Range("M134").Value = -2295.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1116.4166
$ws.Range("I14").Value = 1116.4166
$ws.Range("K14").Value = 3349.2498
$ws.Range("M14").Value = -3176.2498

$ws.Range("H33").Value = 95
$ws.Range("I33").Value = 67.5
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 405
$ws.Range("L33").Value = 900
$ws.Range("M33").Value = -122
$ws.Range("N33").Value = -1466

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H62").Value = 43576.75
$ws.Range("I62").Value = 43576.75
$ws.Range("K62").Value = 43576.75
$ws.Range("M62").Value = -42890.75

$ws.Range("H65").Value = 43576.75
$ws.Range("I65").Value = 43576.75
$ws.Range("K65").Value = 130730.25
$ws.Range("M65").Value = -127298.25

$ws.Range("H70").Value = 8586.799999999999
$ws.Range("I70").Value = 8782.333000000001
$ws.Range("J70").Value = 8456.444
$ws.Range("K70").Value = 8782.333000000001
$ws.Range("L70").Value = 8456.444
$ws.Range("M70").Value = -8512.333000000001
$ws.Range("N70").Value = -8996.444

$ws.Range("H73").Value = 8586.799999999999
$ws.Range("I73").Value = 8782.333000000001
$ws.Range("J73").Value = 8456.444
$ws.Range("K73").Value = 8782.333000000001
$ws.Range("L73").Value = 8456.444
$ws.Range("M73").Value = -7846.333000000001
$ws.Range("N73").Value = -10328.444

$ws.Range("H126").Value = 2999.75
$ws.Range("I126").Value = 2750
$ws.Range("J126").Value = 3249.5
$ws.Range("K126").Value = 8250
$ws.Range("L126").Value = 9748.5
$ws.Range("M126").Value = -5780
$ws.Range("N126").Value = -14688.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1613.5714
$ws.Range("I82").Value = 1613.5714
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1613.5714
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1252.5714
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 1613.5714
$ws.Range("I85").Value = 1613.5714
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1613.5714
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -365.5714
$ws.Range("N85").ClearContents()

$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 1000
$ws.Range("M93").Value = 248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 55000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H62").Value = 14993.6
$ws.Range("I62").Value = 14985
$ws.Range("J62").Value = 14995.75
$ws.Range("K62").Value = 14985
$ws.Range("L62").Value = 14995.75
$ws.Range("M62").Value = -14361
$ws.Range("N62").Value = -16243.75

$ws.Range("H65").Value = 14993.6
$ws.Range("I65").Value = 14985
$ws.Range("J65").Value = 14995.75
$ws.Range("K65").Value = 74925
$ws.Range("L65").Value = 74978.75
$ws.Range("M65").Value = -71805
$ws.Range("N65").Value = -81218.75

$ws.Range("H100").Value = 11113881
$ws.Range("I100").Value = 12501866
$ws.Range("K100").Value = 25003732
$ws.Range("M100").Value = -25003191
